$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.85
$ws.Range("C2").Value = 0.59

$ws.Range("B3").Value = 0.85
$ws.Range("C3").Value = 1.04

$ws.Range("B4").Value = 0.85
$ws.Range("C4").Value = 1.04
